# v1.13: Int. por caución, carga Excel al refrescar, modal marcha de cálculo, G/P acum + interés acum
# Adds one new log row to "Log", one new row to "Resumen" and one new row to "Versiones".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Log sheet: append row 44
# ---------------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("Log")
$wsLog.Range("A44").Value = "27/02/2025"
$wsLog.Range("B44").Value = "21:00"
$wsLog.Range("C44").Value = "Int. por caución y marcha de cálculo"
$wsLog.Range("D44").Value = "Columna Int. por caución en flujo por mes: interés por reinvertir sobrante a un día con tasa de Serie_Cauciones. Carga Excel Serie_Cauciones.xlsx al refrescar (o fallback serie_cauciones.json). Modal al clic en valor mensual con marcha: G/P acum, Int T-1, Base, Tasa, Int T. Cálculo sobre G/P acumulado a la fecha + interés acumulado (reinversión día a día). Fechas ISO (2025-08-25T00:00:00) y columna tasa_diaria."
$wsLog.Range("E44").Value = "Implementacion"

# ---------------------------------------------------------------------------
# Resumen sheet: append row 40
# ---------------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("A40").Value = "Int. por caución"
$wsResumen.Range("B40").Value = "Columna en flujo por mes: interés mensual por colocar el sobrante de caja a la tasa diaria de la serie de cauciones. Carga Serie_Cauciones.xlsx al refrescar (o serie_cauciones.json si no hay Excel). Cálculo: base = G/P acumulado a la fecha + interés acumulado; Int T = base × tasa. Clic en el valor abre modal con marcha (G/P acum, Int T-1, Base, Tasa, Int T)."

# ---------------------------------------------------------------------------
# Versiones sheet: append row 15
# ---------------------------------------------------------------------------
$wsVersiones = $wb.Worksheets.Item("Versiones")
# Leading apostrophe forces text storage (matches existing "1.x" cells, which
# are also stored as text rather than numbers).
$wsVersiones.Range("A15").Value = "'1.13"
$wsVersiones.Range("B15").Value = "27/02/2025"
$wsVersiones.Range("C15").Value = "Int. por caución: columna en flujo, carga Excel al refrescar, modal marcha de cálculo (G/P acum, Base, Tasa, Int T), cálculo sobre G/P acum + interés acum"
